# Insert a new "部门" (department) column in the team_member sheet.
#
# Before:  A=昵称 B=邮箱 C=手机 D=工号 E=职位 F=座位号 G=分机号
# After:   A=昵称 B=邮箱 C=手机 D=工号 E=职位 F=部门  G=座位号 H=分机号
#
# The old "座位号" column (F) keeps its numeric values but slides one
# column to the right into the new "G", and a brand-new "分机号" column
# is appended as "H" carrying what used to live in the old "G". The
# vacated "F" column becomes the new "部门" (department) column with
# per-row text values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- capture existing values for columns F (old 座位号) and G (old 分机号) ---
# before we overwrite anything, since F's old numbers need to move into G,
# and G's old numbers need to move into the brand-new H.
$oldF = @{}
$oldG = @{}
for ($r = 2; $r -le 7; $r++) {
    $oldF[$r] = $ws.Cells.Item($r, 6).Text
    $oldG[$r] = $ws.Cells.Item($r, 7).Text
}

# --- header row ---
$ws.Range("F1").Value = "部门"
$ws.Range("G1").Value = "座位号"
$ws.Range("H1").Value = "分机号"

# --- department values for the new F column ---
$departments = @{
    2 = "管理部"
    3 = "设计部"
    4 = "开发部"
    5 = "客服部"
    6 = "客服部"
    7 = "客服部"
}

for ($r = 2; $r -le 7; $r++) {
    $ws.Cells.Item($r, 6).Value = $departments[$r]   # F: new 部门 text
    $ws.Cells.Item($r, 7).Value = $oldF[$r]           # G: old 座位号 numbers, shifted right
    $ws.Cells.Item($r, 8).Value = $oldG[$r]           # H: old 分机号 numbers, shifted right
}

$ws.Range("F14").Select()
